$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37; existing rows 37-121 shift down to 38-122.
$ws.Rows(37).Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 45133
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 100112013
$ws.Range("G37").Value = "Alcachofa"
$ws.Range("H37").Value = "Madrigal"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = 15000
$ws.Range("N37").Value = "`$/caja 40 unidades"
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 375
$ws.Range("Q37").Value = 40
$ws.Range("R37").Value = "Hortaliza"
